$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A19 to 12.34 and format it as a percentage (0.00%)
$ws.Range("A19").Value = 12.34
$ws.Range("A19").NumberFormat = "0.00%"

# Update selection to A19 (matches the sheetView selection change in the diff)
$ws.Range("A19").Select()
